$wb = $excel.ActiveWorkbook

# Update numeric cells per the source diff (columns H:N = price/profit calc columns).
# For cells that are removed entirely in the diff (no longer present in the row),
# we use ClearContents() so the cell reverts to empty rather than holding a stale 0.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3490
$ws.Range("J17").Value = 3490
$ws.Range("L17").Value = 10470
$ws.Range("N17").Value = -10806

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H92").Value = 690.2857
$ws.Range("J92").Value = 377.5
$ws.Range("L92").Value = 377.5
$ws.Range("N92").Value = -2873.5

$ws.Range("H118").Value = 250
$ws.Range("I118").Value = 250
$ws.Range("K118").Value = 750
$ws.Range("M118").Value = 907

$ws.Range("H129").Value = 2238.125
$ws.Range("J129").Value = 2405.8333
$ws.Range("L129").Value = 7217.499899999999
$ws.Range("N129").Value = -17217.4999

$ws.Range("H132").Value = 1643.5
$ws.Range("I132").Value = 1180.625
$ws.Range("K132").Value = 3541.875
$ws.Range("M132").Value = -1011.875

$ws.Range("H137").Value = 2031.421
$ws.Range("I137").Value = 2093.5625
$ws.Range("K137").Value = 6280.6875
$ws.Range("M137").Value = -3730.6875

$ws.Range("H138").Value = 2468.2632
$ws.Range("J138").Value = 3776
$ws.Range("L138").Value = 11328
$ws.Range("N138").Value = -21608

$ws.Range("H141").Value = 1348.5264
$ws.Range("I141").Value = 1348.5264
$ws.Range("K141").Value = 4045.5792
$ws.Range("M141").Value = 1134.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1963.907
$ws.Range("I32").Value = 1848.7
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 1848.7
$ws.Range("L32").Value = 3500
$ws.Range("M32").Value = -1561.7
$ws.Range("N32").Value = -4074

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H107").Value = 4976.25
$ws.Range("I107").Value = 4911
$ws.Range("J107").Value = 4998
$ws.Range("K107").Value = 4911
$ws.Range("L107").Value = 4998
$ws.Range("M107").Value = -2991
$ws.Range("N107").Value = -8838

$ws.Range("H134").Value = 2906.6667
$ws.Range("I134").Value = 1427.1428
$ws.Range("J134").Value = 4978
$ws.Range("K134").Value = 4281.428400000001
$ws.Range("L134").Value = 14934
$ws.Range("M134").Value = -1746.428400000001
$ws.Range("N134").Value = -20004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2695
$ws.Range("I16").Value = 2495.6667
$ws.Range("K16").Value = 2495.6667
$ws.Range("M16").Value = -2208.6667

$ws.Range("H31").Value = 2425.4583
$ws.Range("I31").Value = 1138.375
$ws.Range("K31").Value = 1138.375
$ws.Range("M31").Value = -843.375

$ws.Range("H34").Value = 2425.4583
$ws.Range("I34").Value = 1138.375
$ws.Range("K34").Value = 1138.375
$ws.Range("M34").Value = -936.375

$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797

$ws.Range("H113").Value = 2695
$ws.Range("I113").Value = 2495.6667
$ws.Range("K113").Value = 2495.6667
$ws.Range("M113").Value = -325.6667000000002

$ws.Range("H132").Value = 3663.5
$ws.Range("I132").Value = 3663.5
$ws.Range("K132").Value = 10990.5
$ws.Range("M132").Value = -8460.5

$ws.Range("H134").Value = 1566.6666
$ws.Range("I134").Value = 1850
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 5550
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -3015
$ws.Range("N134").Value = -8070

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 106
$ws.Range("I2").Value = 101
$ws.Range("J2").Value = 106.07692
$ws.Range("K2").Value = 606
$ws.Range("L2").Value = 636.4615200000001
$ws.Range("M2").Value = -493
$ws.Range("N2").Value = -862.4615200000001

$ws.Range("H51").Value = 250
$ws.Range("I51").Value = 250
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 750
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -290
$ws.Range("N51").ClearContents()

$ws.Range("H122").Value = 2171.8333
$ws.Range("I122").Value = 416.3
$ws.Range("J122").Value = 10949.5
$ws.Range("K122").Value = 3746.7
$ws.Range("L122").Value = 98545.5
$ws.Range("M122").Value = -1296.7
$ws.Range("N122").Value = -103445.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 30767
$ws.Range("J128").Value = 30767
$ws.Range("L128").Value = 30767
$ws.Range("N128").Value = -40727

$ws.Range("H132").Value = 3224.5
$ws.Range("I132").Value = 3224.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9673.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7143.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2950.375
$ws.Range("I40").Value = 3086.1428
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 3086.1428
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2950.1428
$ws.Range("N40").Value = -2272

$ws.Range("H132").Value = 1127.125
$ws.Range("I132").Value = 1028.1666
$ws.Range("J132").Value = 1424
$ws.Range("K132").Value = 3084.4998
$ws.Range("L132").Value = 4272
$ws.Range("M132").Value = -554.4998000000001
$ws.Range("N132").Value = -9332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10336

$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30298

$ws.Range("H132").Value = 1839.2812
$ws.Range("I132").Value = 960.9583
$ws.Range("J132").Value = 4474.25
$ws.Range("K132").Value = 2882.8749
$ws.Range("L132").Value = 13422.75
$ws.Range("M132").Value = -352.8748999999998
$ws.Range("N132").Value = -18482.75

$ws.Range("H136").Value = 1017.35
$ws.Range("I136").Value = 1017.35
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3052.05
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -502.0500000000002
$ws.Range("N136").ClearContents()
